$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Activation date: "01/01/2012" -> "01/01/2023" (affects B8/C8 and B13/C13,
#    which share the same text). Assigning a date-shaped string straight to
#    .Value makes Excel coerce it to a date serial, so instead we stage a
#    text-producing formula in a scratch cell and copy only its *value* onto
#    the target cells - this keeps the original General number format /
#    style (s=2 / s=3) on B8/C8/B13/C13 intact while storing a real text
#    value (no new number formats get added to styles.xml).
# ---------------------------------------------------------------------------
$scratch = $ws.Cells.Item(100, 26)  # Z100, well outside the used range
$scratch.Formula = "=""01/01/2023"""
$scratch.Copy()

$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").PasteSpecial(-4163)

$scratch.Clear()

# ---------------------------------------------------------------------------
# 2) New English "Objectives" text under row 11 (A11 = "Objectives:")
# ---------------------------------------------------------------------------
$objectives = "Provide the student with the basic knowledge of magnetic and superconducting materials aiming their application in devices."
$ws.Cells.Item(11, 2).Value = $objectives
$ws.Cells.Item(11, 3).Value = $objectives
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) New English "Short syllabus" text under row 14 (A14 = "Short syllabus:")
# ---------------------------------------------------------------------------
$syllabus = "Fundamental concepts of magnetic properties of matter. Electron magnetism. Ferromagnetism. Magnetic Materials and Applications: soft and hard. Exchange interaction in oxides and metals. Magnetism - Classical Phenomenology: diamagnetism and paramagnetism. Magnetism - Quantum Phenomenology: ferromagnetism. Magnetic Anisotropy and Spin-Orbit Interaction. Magnetostriction and magnetostrictive materials -Introduction and applications. Basic concepts of superconductivity. Superconductivity - Quantum Origin. Superwave – Consequences. Quantum Interference – SQUID. Superconducting Materials and Applications"
$ws.Cells.Item(14, 2).Value = $syllabus
$ws.Cells.Item(14, 3).Value = $syllabus
$ws.Range("B15:C15").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 4) Same text reused under row 16 (A16 = "Syllabus:")
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 2).Value = $syllabus
$ws.Cells.Item(16, 3).Value = $syllabus
$ws.Range("B15:C15").Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 5) Updated "Norma de recuperacao" text in row 20
# ---------------------------------------------------------------------------
$recuperacao = "A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável."
$ws.Cells.Item(20, 2).Value = $recuperacao
$ws.Cells.Item(20, 3).Value = $recuperacao

$excel.CutCopyMode = 0
